$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '39.157.51'
$ws.Range("E2").Value = '  -2.51%  '
$ws.Range("D3").Value = '2.193.64'
$ws.Range("E3").Value = '  -5.99%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '293.80'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -4.55%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '80.64'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -5.49%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.508'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -4.17%  '
$ws.Range("E8").Value = '  -0.09%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.464'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -4.19%  '
$ws.Range("E10").Value = '  -6.43%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '28.90'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -4.19%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '46.83'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -11.07%  '
$ws.Range("E13").Value = '  -2.77%  '
$ws.Range("D14").Value = '2.534.45'
$ws.Range("E14").Value = '  -5.81%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.18'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.96%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '13.84'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -6.06%  '
$ws.Range("D17").Value = '2.193.66'
$ws.Range("E17").Value = '  -5.89%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.707'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -6.26%  '
$ws.Range("D19").Value = '39.057.13'
$ws.Range("E19").Value = '  -2.58%  '
$ws.Range("D20").Value = '0.0₃0866'
$ws.Range("E20").Value = '  -4.33%  '
$ws.Range("E21").Value = '  -6.89%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '64.42'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.81%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.21'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -4.55%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '224.74'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -4.85%  '
$ws.Range("E25").Value = '  -0.14%  '
$ws.Range("E26").Value = '  -7.26%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.79'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.01%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '22.37'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.47%  '
$ws.Range("E29").Value = '  -2.14%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.98'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.08%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '148.72'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.44%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '31.47'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -10.62%  '
$ws.Range("E33").Value = '  -0.10%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.76'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -7.14%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.34'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.33%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0690'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -4.90%  '
$ws.Range("E37").Value = '  -4.00%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '15.14'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.09%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0949'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -5.04%  '
$ws.Range("E40").Value = '  -5.34%  '
$ws.Range("E41").Value = '  -4.35%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.58'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -5.99%  '
$ws.Range("D43").Value = '1.894.97'
$ws.Range("E43").Value = '  -2.20%  '
$ws.Range("E44").Value = '  -9.07%  '
$ws.Range("E45").Value = '  -3.72%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.97'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.25%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '15.83'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -11.22%  '
$ws.Range("E48").Value = '  -3.90%  '
$ws.Range("D49").Value = '2.408.24'
$ws.Range("E49").Value = '  -5.90%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '70.81'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.92%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '86.77'
$ws.Range("D51").Style = "Normal"
